$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-10 (columns A-T).
# Sending/Target clusters now span ECs, FAPs, sCs (previously only FAPs, sCs),
# so the 2x2 grid of sending/target combinations grows to a 3x3 grid.
$rows = @(
  @("ECs","Angpt1","Tek","ECs",2,0.6666666666666666,0.1659303333333333,0.497791,0.01127931880021135,0.01127931880021135,3,1,34.10446833333334,102.313405,0.5118942073015388,0.5118942073015389,5.658965798706111,50.930692188355,0.005773817956135533,0.005773817956135535),
  @("ECs","Angpt1","Tek","FAPs",2,0.6666666666666666,0.1659303333333333,0.497791,0.01127931880021135,0.01127931880021135,3,1,31.083557,93.25067100000001,0.4665515560925921,0.4665515560925922,5.157704974195666,46.419344767761,0.005262383737903035,0.005262383737903037),
  @("ECs","Angpt1","Tek","sCs",2,0.6666666666666666,0.1659303333333333,0.497791,0.01127931880021135,0.01127931880021135,3,1,1.436030666666667,4.308092,0.02155423660586901,0.02155423660586901,0.2382810471968889,2.144529424772,0.000243117106172782,0.000243117106172782),
  @("FAPs","Angpt1","Tek","ECs",3,1,12.95459633333333,38.863789,0.8806046431437031,0.8806046431437031,3,1,34.10446833333334,102.313405,0.5118942073015388,0.5118942073015389,441.8096204212828,3976.286583791545,0.4507764157481003,0.4507764157481005),
  @("FAPs","Angpt1","Tek","FAPs",3,1,12.95459633333333,38.863789,0.8806046431437031,0.8806046431437031,3,1,31.083557,93.25067100000001,0.4665515560925921,0.4665515560925922,402.6749335391577,3624.074401852419,0.4108474665610565,0.4108474665610566),
  @("FAPs","Angpt1","Tek","sCs",3,1,12.95459633333333,38.863789,0.8806046431437031,0.8806046431437031,3,1,1.436030666666667,4.308092,0.02155423660586901,0.02155423660586901,18.60319760895422,167.428778480588,0.01898076083454622,0.01898076083454622),
  @("sCs","Angpt1","Tek","ECs",3,1,1.590497666666667,4.771493,0.1081160380560855,0.1081160380560855,3,1,34.10446833333334,102.313405,0.5118942073015388,0.5118942073015389,54.24307730707389,488.187695763665,0.05534397359730289,0.05534397359730291),
  @("sCs","Angpt1","Tek","FAPs",3,1,1.590497666666667,4.771493,0.1081160380560855,0.1081160380560855,3,1,31.083557,93.25067100000001,0.4665515560925921,0.4665515560925922,49.43832488020033,444.944923921803,0.05044170579363261,0.05044170579363261),
  @("sCs","Angpt1","Tek","sCs",3,1,1.590497666666667,4.771493,0.1081160380560855,0.1081160380560855,3,1,1.436030666666667,4.308092,0.02155423660586901,0.02155423660586901,2.284003424595111,20.556030821356,0.002330358665150005,0.002330358665150005)
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

Write-Host "Wrote $($rows.Length) data rows"
